$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- swap row 26 and row 27 (all columns except A) ---
$ws.Cells.Item(26, 2).Value = 6732711
$ws.Cells.Item(26, 6).Value = "Banga Gargzdai"
$ws.Cells.Item(26, 7).Value = "FK Zalgiris Vilnius"
$ws.Cells.Item(26, 8).Value = 1
$ws.Cells.Item(26, 9).Value = 4
$ws.Cells.Item(26, 12).Value = 3.6
$ws.Cells.Item(26, 13).Value = 1.571
$ws.Cells.Item(26, 14).Value = 11
$ws.Cells.Item(26, 15).Value = 4.75
$ws.Cells.Item(26, 16).Value = 1.25
$ws.Cells.Item(26, 17).Value = 1.5
$ws.Cells.Item(26, 18).Value = 1.975
$ws.Cells.Item(26, 19).Value = 1.825
$ws.Cells.Item(26, 21).Value = 1.8
$ws.Cells.Item(26, 22).Value = 2
$ws.Cells.Item(26, 25).Value = 0.25
$ws.Cells.Item(26, 26).Value = -1
$ws.Cells.Item(26, 27).Value = 0.825
$ws.Cells.Item(26, 28).Value = 0.8
$ws.Cells.Item(26, 29).Value = -1
$ws.Cells.Item(27, 2).Value = 6732773
$ws.Cells.Item(27, 6).Value = "Suduva Marijampole"
$ws.Cells.Item(27, 7).Value = "Hegelmann Litauen"
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 9).Value = 1
$ws.Cells.Item(27, 12).Value = 3.8
$ws.Cells.Item(27, 13).Value = 1.533
$ws.Cells.Item(27, 14).Value = 5
$ws.Cells.Item(27, 15).Value = 4.2
$ws.Cells.Item(27, 16).Value = 1.533
$ws.Cells.Item(27, 17).Value = 1
$ws.Cells.Item(27, 18).Value = 1.875
$ws.Cells.Item(27, 19).Value = 1.925
$ws.Cells.Item(27, 21).Value = 1.9
$ws.Cells.Item(27, 22).Value = 1.9
$ws.Cells.Item(27, 25).Value = 0.5329999999999999
$ws.Cells.Item(27, 26).Value = 0
$ws.Cells.Item(27, 27).Value = -0
$ws.Cells.Item(27, 28).Value = -1
$ws.Cells.Item(27, 29).Value = 0.8999999999999999

# --- swap row 89 and row 90 (all columns except A) ---
$ws.Cells.Item(89, 2).Value = 7326568
$ws.Cells.Item(89, 6).Value = "Hegelmann Litauen"
$ws.Cells.Item(89, 7).Value = "Panevezys"
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = "D"
$ws.Cells.Item(89, 11).Value = 2.375
$ws.Cells.Item(89, 12).Value = 3.2
$ws.Cells.Item(89, 13).Value = 2.625
$ws.Cells.Item(89, 14).Value = 2.7
$ws.Cells.Item(89, 15).Value = 3.2
$ws.Cells.Item(89, 16).Value = 2.3
$ws.Cells.Item(89, 17).Value = 0
$ws.Cells.Item(89, 18).Value = 2.05
$ws.Cells.Item(89, 19).Value = 1.75
$ws.Cells.Item(89, 20).Value = 2.25
$ws.Cells.Item(89, 21).Value = 1.875
$ws.Cells.Item(89, 22).Value = 1.925
$ws.Cells.Item(89, 24).Value = 2.2
$ws.Cells.Item(89, 25).Value = -1
$ws.Cells.Item(89, 26).Value = 0
$ws.Cells.Item(89, 27).Value = -0
$ws.Cells.Item(89, 29).Value = 0.925
$ws.Cells.Item(90, 2).Value = 6732827
$ws.Cells.Item(90, 6).Value = "FK Dziugas Telsiai"
$ws.Cells.Item(90, 7).Value = "FK Kauno Zalgiris"
$ws.Cells.Item(90, 9).Value = 2
$ws.Cells.Item(90, 10).Value = "A"
$ws.Cells.Item(90, 11).Value = 6
$ws.Cells.Item(90, 12).Value = 3.9
$ws.Cells.Item(90, 13).Value = 1.444
$ws.Cells.Item(90, 14).Value = 4.75
$ws.Cells.Item(90, 15).Value = 3.6
$ws.Cells.Item(90, 16).Value = 1.65
$ws.Cells.Item(90, 17).Value = 0.75
$ws.Cells.Item(90, 18).Value = 1.9
$ws.Cells.Item(90, 19).Value = 1.9
$ws.Cells.Item(90, 20).Value = 2.5
$ws.Cells.Item(90, 21).Value = 1.95
$ws.Cells.Item(90, 22).Value = 1.85
$ws.Cells.Item(90, 24).Value = -1
$ws.Cells.Item(90, 25).Value = 0.6499999999999999
$ws.Cells.Item(90, 26).Value = -1
$ws.Cells.Item(90, 27).Value = 0.8999999999999999
$ws.Cells.Item(90, 29).Value = 0.8500000000000001

# --- swap row 103 and row 104 (all columns except A) ---
$ws.Cells.Item(103, 2).Value = 6732837
$ws.Cells.Item(103, 6).Value = "Suduva Marijampole"
$ws.Cells.Item(103, 7).Value = "FK Riteriai"
$ws.Cells.Item(103, 8).Value = 0
$ws.Cells.Item(103, 9).Value = 3
$ws.Cells.Item(103, 10).Value = "A"
$ws.Cells.Item(103, 11).Value = 3.6
$ws.Cells.Item(103, 12).Value = 3.6
$ws.Cells.Item(103, 13).Value = 1.8
$ws.Cells.Item(103, 14).Value = 3
$ws.Cells.Item(103, 15).Value = 3.6
$ws.Cells.Item(103, 16).Value = 2
$ws.Cells.Item(103, 18).Value = 2
$ws.Cells.Item(103, 19).Value = 1.8
$ws.Cells.Item(103, 20).Value = 2.5
$ws.Cells.Item(103, 21).Value = 1.975
$ws.Cells.Item(103, 22).Value = 1.825
$ws.Cells.Item(103, 23).Value = -1
$ws.Cells.Item(103, 25).Value = 1
$ws.Cells.Item(103, 26).Value = -1
$ws.Cells.Item(103, 27).Value = 0.8
$ws.Cells.Item(103, 28).Value = 0.9750000000000001
$ws.Cells.Item(104, 2).Value = 7465686
$ws.Cells.Item(104, 6).Value = "FK Kauno Zalgiris"
$ws.Cells.Item(104, 7).Value = "Hegelmann Litauen"
$ws.Cells.Item(104, 8).Value = 4
$ws.Cells.Item(104, 9).Value = 2
$ws.Cells.Item(104, 10).Value = "H"
$ws.Cells.Item(104, 11).Value = 2.3
$ws.Cells.Item(104, 12).Value = 4
$ws.Cells.Item(104, 13).Value = 2.3
$ws.Cells.Item(104, 14).Value = 2.55
$ws.Cells.Item(104, 15).Value = 4
$ws.Cells.Item(104, 16).Value = 2.2
$ws.Cells.Item(104, 18).Value = 1.8
$ws.Cells.Item(104, 19).Value = 2
$ws.Cells.Item(104, 20).Value = 2.75
$ws.Cells.Item(104, 21).Value = 1.85
$ws.Cells.Item(104, 22).Value = 1.95
$ws.Cells.Item(104, 23).Value = 1.55
$ws.Cells.Item(104, 25).Value = -1
$ws.Cells.Item(104, 26).Value = 0.8
$ws.Cells.Item(104, 27).Value = -1
$ws.Cells.Item(104, 28).Value = 0.8500000000000001

# --- add new row 128 (copy style from row 127 for columns A and E) ---
$ws.Cells.Item(127, 1).Copy() | Out-Null
$ws.Cells.Item(128, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(127, 5).Copy() | Out-Null
$ws.Cells.Item(128, 5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(128, 1).Value = 126
$ws.Cells.Item(128, 2).Value = 7862040
$ws.Cells.Item(128, 3).Value = "Lithuania A Lyga"
$ws.Cells.Item(128, 4).Value = "Lithuania A Lyga"
$ws.Cells.Item(128, 5).Value = 45382.375
$ws.Cells.Item(128, 6).Value = "Suduva Marijampole"
$ws.Cells.Item(128, 7).Value = "FK Dainava Alytus"
$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(128, 9).Value = 1
$ws.Cells.Item(128, 10).Value = "A"
$ws.Cells.Item(128, 11).Value = 2.45
$ws.Cells.Item(128, 12).Value = 3.25
$ws.Cells.Item(128, 13).Value = 2.5
$ws.Cells.Item(128, 14).Value = 2.15
$ws.Cells.Item(128, 15).Value = 3.2
$ws.Cells.Item(128, 16).Value = 3
$ws.Cells.Item(128, 17).Value = -0.25
$ws.Cells.Item(128, 18).Value = 1.9
$ws.Cells.Item(128, 19).Value = 1.9
$ws.Cells.Item(128, 20).Value = 2
$ws.Cells.Item(128, 21).Value = 1.85
$ws.Cells.Item(128, 22).Value = 1.95
$ws.Cells.Item(128, 23).Value = -1
$ws.Cells.Item(128, 24).Value = -1
$ws.Cells.Item(128, 25).Value = 2
$ws.Cells.Item(128, 26).Value = -1
$ws.Cells.Item(128, 27).Value = 0.8999999999999999
$ws.Cells.Item(128, 28).Value = -1
$ws.Cells.Item(128, 29).Value = 0.95

# --- add new row 129 (copy style from row 127 for columns A and E) ---
$ws.Cells.Item(127, 1).Copy() | Out-Null
$ws.Cells.Item(129, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(127, 5).Copy() | Out-Null
$ws.Cells.Item(129, 5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(129, 1).Value = 127
$ws.Cells.Item(129, 2).Value = 7862042
$ws.Cells.Item(129, 3).Value = "Lithuania A Lyga"
$ws.Cells.Item(129, 4).Value = "Lithuania A Lyga"
$ws.Cells.Item(129, 5).Value = 45382.51736111111
$ws.Cells.Item(129, 6).Value = "FK Kauno Zalgiris"
$ws.Cells.Item(129, 7).Value = "FK Zalgiris Vilnius"
$ws.Cells.Item(129, 8).Value = 0
$ws.Cells.Item(129, 9).Value = 1
$ws.Cells.Item(129, 10).Value = "A"
$ws.Cells.Item(129, 11).Value = 4
$ws.Cells.Item(129, 12).Value = 3.5
$ws.Cells.Item(129, 13).Value = 1.727
$ws.Cells.Item(129, 14).Value = 3.3
$ws.Cells.Item(129, 15).Value = 3.2
$ws.Cells.Item(129, 16).Value = 2
$ws.Cells.Item(129, 17).Value = 0.25
$ws.Cells.Item(129, 18).Value = 1.975
$ws.Cells.Item(129, 19).Value = 1.825
$ws.Cells.Item(129, 20).Value = 2.5
$ws.Cells.Item(129, 21).Value = 1.9
$ws.Cells.Item(129, 22).Value = 1.9
$ws.Cells.Item(129, 23).Value = -1
$ws.Cells.Item(129, 24).Value = -1
$ws.Cells.Item(129, 25).Value = 1
$ws.Cells.Item(129, 26).Value = -1
$ws.Cells.Item(129, 27).Value = 0.825
$ws.Cells.Item(129, 28).Value = -1
$ws.Cells.Item(129, 29).Value = 0.8999999999999999
